$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9976536631584167
$ws.Range("B1").Value = 2.834280967712402
$ws.Range("C1").Value = 3.059633016586304
$ws.Range("D1").Value = 3.6593177318573
$ws.Range("E1").Value = 1.404904365539551
